# netCrypto.xlsx update:
#  - SheetName1!T2 value changes from 478002 to 501242
#  - The active selection on the sheet moves from T2 to T3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 501242
$ws.Range("T3").Select()
